# Change the sample row from the generic "degiro" user to the
# "hyves" employer (username, email, name, address, postal code).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "hyves"
$ws.Range("B2").Value = "info@hyves.nl"
$ws.Range("C2").Value = "Hyves"
$ws.Range("D2").Value = "Basisweg 30"
$ws.Range("E2").Value = "1043AP"
# F2 (city, "Amsterdam") is unchanged.

# The mailto: hyperlink on B2 needs to point at the new address too.
# Setting properties on the existing Hyperlinks item only appends a
# duplicate relationship, so drop every hyperlink on the sheet first and
# re-add the single link we need.
$ws.Hyperlinks.Delete()
[void]$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:info@hyves.nl", [System.Type]::Missing, [System.Type]::Missing, "info@hyves.nl")

# D3 previously carried its own one-off font; bring it back in line with
# the rest of column D (wrap text, shared default font) like D4:D16.
$ws.Range("D3").WrapText = $true

# Move the active selection to D12.
[void]$ws.Range("D12").Select()
